# CIERRE 6 JUL 22
# Update the incentive-payment voucher sheet for the new closing period
# (MAYO -> JUNIO), move the active tab from "ARQUITECTO" to
# "VALES DE INSENTIVOS", and move the selection there to A6.

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# Update the month text on the "VALES DE INSENTIVOS" sheet.
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE  JUNIO   2022"

# Switch the active sheet/tab to "VALES DE INSENTIVOS" and move the
# selection there to A6 (was previously on "ARQUITECTO", selection A11:B11).
$wsVales.Activate() | Out-Null
$wsVales.Range("A6").Select() | Out-Null
